# CalSim3DataExtractionInitFile_v4.xlsx - "updated init files up to scenario 63"
#
# The scenario-listing / variable-block cell references in column D
# (rows 5-11) point at the block of cells ("Upper Left Cell" .. the row-31
# anchors of the Scenario Listings tab) that holds the scenario index /
# name / directory / DSS-path / date lists. Bumping the run up to
# "scenario 63" moved that reference block down two rows (31 -> 33) in the
# external listing workbook, so each of the seven anchor cells here is
# updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

$ws.Range("D5").Value  = "A33"
$ws.Range("D6").Value  = "B33"
$ws.Range("D7").Value  = "C33"
$ws.Range("D8").Value  = "G33"
$ws.Range("D9").Value  = "H33"
$ws.Range("D10").Value = "I33"
$ws.Range("D11").Value = "J33"

# Restore/keep the original selection block over the edited cells and
# scroll the view down a few rows, matching the editor's on-screen state
# when the file was saved.
$ws.Activate() | Out-Null
$ws.Range("D5:D11").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
